$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from the last existing data row (181) down into the
# two new rows (182, 183) so the new cells pick up the same cell styles
# (date format in column A, default body style in B:AH) as the rest of the table.
$ws.Range("A181:AH181").Copy()
$ws.Range("A182:AH183").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 182
$ws.Range("A182").Value = 44957
$ws.Range("B182").Value = 0.091714199939958
$ws.Range("C182").Value = 0.06904603497126827
$ws.Range("D182").Value = 0.05537432066545911
$ws.Range("E182").Value = 0.05537459283387625
$ws.Range("F182").Value = 0.04034715420611601
$ws.Range("G182").Value = 0.06175283240005203
$ws.Range("H182").Value = 0.0710132478773664
$ws.Range("I182").Value = 0.07299916055520494
$ws.Range("J182").Value = 0.09691945361408894
$ws.Range("K182").Value = 0.06786610708816965
$ws.Range("L182").Value = 0.08051215835961179
$ws.Range("M182").Value = 0.0785461845709865
$ws.Range("N182").Value = 0.02565695213959773
$ws.Range("O182").Value = 0.03909594201244215
$ws.Range("P182").Value = 0
$ws.Range("Q182").Value = 0.0003078404597083306
$ws.Range("R182").Value = 0.003463904694816611
$ws.Range("S182").Value = -0.01374443666798697
$ws.Range("T182").Value = 0
$ws.Range("U182").Value = 0.09949496313500972
$ws.Range("V182").Value = 0.058805
$ws.Range("W182").Value = 0.109055
$ws.Range("X182").Value = -0.018539
$ws.Range("Y182").Value = 0.023481
$ws.Range("Z182").Value = 0.005496999999999999
$ws.Range("AA182").Value = 0.006319
$ws.Range("AB182").Value = 0.001666666666666667
$ws.Range("AC182").Value = -0.01566591597069998
$ws.Range("AD182").Value = 0
$ws.Range("AE182").Value = 0.0468490262906176
$ws.Range("AF182").Value = 0.05698371271151981
$ws.Range("AG182").Value = 0.06837699999999999
$ws.Range("AH182").Value = 0.09091100000000001

# Row 183
$ws.Range("A183").Value = 44985
$ws.Range("B183").Value = -0.06256015399422521
$ws.Range("C183").Value = -0.05448975178412663
$ws.Range("D183").Value = -0.05137564916524462
$ws.Range("E183").Value = -0.05070546737213399
$ws.Range("F183").Value = -0.04144982870890968
$ws.Range("G183").Value = -0.0261124466467153
$ws.Range("H183").Value = -0.02984407051529225
$ws.Range("I183").Value = -0.02894033027879028
$ws.Range("J183").Value = -0.0180926475650186
$ws.Range("K183").Value = -0.02501167842393226
$ws.Range("L183").Value = -0.02225724133990969
$ws.Range("M183").Value = -0.06542898691226373
$ws.Range("N183").Value = 0.00629740755803132
$ws.Range("O183").Value = -0.0129352960562451
$ws.Range("P183").Value = 0
$ws.Range("Q183").Value = 0.004955495233886831
$ws.Range("R183").Value = 0.000817719055509869
$ws.Range("S183").Value = 0.01856550054132899
$ws.Range("T183").Value = 0
$ws.Range("U183").Value = -0.06037680233843445
$ws.Range("V183").Value = -0.027521
$ws.Range("W183").Value = -0.08667799999999999
$ws.Range("X183").Value = 0.013952
$ws.Range("Y183").Value = -0.000687
$ws.Range("Z183").Value = -0.000316
$ws.Range("AA183").Value = -0.02187
$ws.Range("AB183").Value = 0.001666666666666667
$ws.Range("AC183").Value = 0.005445778553076633
$ws.Range("AD183").Value = 0
$ws.Range("AE183").Value = 0.0363027081212588
$ws.Range("AF183").Value = 0.03414226202012283
$ws.Range("AG183").Value = -0.052326
$ws.Range("AH183").Value = -0.062184

# Extend the conditional-formatting ranges to cover the two newly added rows,
# mirroring how Excel keeps a "live" formatted range in sync with new data.
$dateCf = $ws.Range("A1:A181").FormatConditions.Item(1)
$dateCf.ModifyAppliesToRange($ws.Range("A1:A183"))

$bodyCf = $ws.Range("B2:AH181").FormatConditions.Item(1)
$bodyCf.ModifyAppliesToRange($ws.Range("B2:AH183"))

